# Reorder columns and fix uwls3 values in the stats table on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order (row 1), and the corresponding row 2 values.
$headers = @("fishers_z_est", "fishers_z_t_value", "hsma_est", "hsma_t_value", "re_est", "re_t_value", "uwls_est", "uwls_t_value", "uwls3_est", "uwls3_t_value")
$values  = @(0.06030813277927358, 16.24184553868096, 0.04567022432261091, 17.6906469956982, 0.05914186723756871, 16.42125501540017, 0.04887186309074776, 17.52923840917323, 0.04872269871584915, 17.56199708184831)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}
